# Update yearly.xlsx "Overview" sheet:
#   - roll the 5-year reporting window forward by one fiscal year
#     (drop 1396/12, shift 1397/12..1400/12 left, append 1401/12)
#   - refresh each metric's 5 trailing-year values to match the new window,
#     appending the newly read figure for 1401/12 in column I

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column headers (year labels) for both tables on the sheet ----
$ws.Range("E8").Value2  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value2  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value2  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value2  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value2  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value2 = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value2 = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value2 = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value2 = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value2 = "دوازده ماهه منتهی به 1401/12"

# ---- هزینه حمل و نقل و انتقال (row 10) ----
$ws.Range("E10").Value2 = 2689
$ws.Range("F10").Value2 = 447
$ws.Range("G10").Value2 = 17174
$ws.Range("H10").Value2 = 15466
$ws.Range("I10").Value2 = 13950

# ---- حق العمل و کمیسیون فروش (row 13) ----
$ws.Range("E13").Value2 = 413
$ws.Range("F13").Value2 = 1467
$ws.Range("G13").Value2 = 5436
$ws.Range("H13").Value2 = 5562
$ws.Range("I13").Value2 = 8044

# ---- هزینه انرژی (آب، برق، گاز و سوخت) (row 15) ----
$ws.Range("E15").Value2 = 0
$ws.Range("F15").Value2 = 0
$ws.Range("G15").Value2 = 353
$ws.Range("H15").Value2 = 378
$ws.Range("I15").Value2 = 662

# ---- هزینه استهلاک (row 16) ----
$ws.Range("E16").Value2 = 1778
$ws.Range("F16").Value2 = 1975
$ws.Range("G16").Value2 = 2177
$ws.Range("H16").Value2 = 2408
$ws.Range("I16").Value2 = 2355

# ---- هزینه حقوق و دستمزد (row 17) ----
$ws.Range("E17").Value2 = 35986
$ws.Range("F17").Value2 = 45487
$ws.Range("G17").Value2 = 62908
$ws.Range("H17").Value2 = 94267
$ws.Range("I17").Value2 = 136021

# ---- سایر هزینه ها (row 19) ----
$ws.Range("E19").Value2 = 15190
$ws.Range("F19").Value2 = 18570
$ws.Range("G19").Value2 = 29181
$ws.Range("H19").Value2 = 45766
$ws.Range("I19").Value2 = 54852

# ---- جمع (row 20) ----
$ws.Range("E20").Value2 = 56056
$ws.Range("F20").Value2 = 67946
$ws.Range("G20").Value2 = 117229
$ws.Range("H20").Value2 = 163847
$ws.Range("I20").Value2 = 215884

# ---- تعداد پرسنل غیر تولیدی شرکت (row 26) ----
$ws.Range("E26").Value2 = 135
$ws.Range("F26").Value2 = 129
$ws.Range("G26").Value2 = 124
$ws.Range("H26").Value2 = 126
$ws.Range("I26").Value2 = 123

# ---- تعداد پرسنل تولیدی شرکت (row 27) ----
$ws.Range("E27").Value2 = 61
$ws.Range("F27").Value2 = 61
$ws.Range("G27").Value2 = 62
$ws.Range("H27").Value2 = 62
$ws.Range("I27").Value2 = 60
